$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new header "Save" in H1, matching the formatting of the existing header cells (e.g. G1)
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add the data values for the new "Save" column
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
